# Bug fixes to get it running again
#
# - Rename "Sheet1" to "Configuration Testing" (it's the active/first sheet)
# - Fill in the procedure/result notes for rows 2-3 (Test/Procedure columns)
# - Grow rows 2-3 to fit the new wrapped text
# - Leave the active selection on B3, matching where the author ended up editing

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Sheet1")
$ws.Name = "Configuration Testing"

$ws.Range("B2").Value = "Go into the LbcbPlugin Folder and double click on the LbcbPlugin code file"
$ws.Range("C2").Value = "MATLAB should start up with a command window and an editor window"
$ws.Range("C3").Value = "The window should clear."
$ws.Range("B3").Value = "You can dismiss the editor window.  Type ""clearSpace"" in the command window.  "

$ws.Rows.Item(2).RowHeight = 30
$ws.Rows.Item(3).RowHeight = 30

$ws.Range("B3").Select()
